$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.772.77'
$ws.Range("E2").Value = '  +0.81%  '
$ws.Range("D3").Value = '1.848.87'
$ws.Range("E3").Value = '  +0.20%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.007'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.06%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '335.53'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.37%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.006'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.01%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4669'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +0.78%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3852'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -0.27%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '46.93'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +2.00%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.07914'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +0.35%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.9673'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -3.13%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '21.29'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -1.00%  '
$ws.Range("D13").Value = '1.867.65'
$ws.Range("E13").Value = '  +1.40%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.866'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -1.70%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.124'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -0.12%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '1.007'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -0.14%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '90.82'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +2.68%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.06617'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -0.99%  '
$ws.Range("E19").Value = '  -0.70%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '17.25'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +0.63%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '1.007'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -0.01%  '
$ws.Range("D22").Value = '27.765.20'
$ws.Range("E22").Value = '  +0.79%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.337'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.97%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '10.79'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -0.91%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.291'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.73%  '
$ws.Range("D26").Value = '2.085.46'
$ws.Range("E26").Value = '  +1.22%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '159.24'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.28%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '19.44'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -0.38%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.064'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -2.52%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '5.376'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -0.80%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '118.53'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -1.01%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.09433'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +0.25%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.9384'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -3.98%  '
$ws.Range("E34").Value = '  -0.07%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '5.251'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -1.10%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.326'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -0.27%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.06017'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -0.15%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.02212'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -0.73%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '8.207'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -1.16%  '
$ws.Range("E40").Value = '  +0.06%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.161'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -1.58%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.5799'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -1.80%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.1845'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -0.97%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '10.06'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -2.88%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.285'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +3.70%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '11.93'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -1.60%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.5443'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -2.55%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.931'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +1.06%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.06837'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +2.05%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '110.77'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +0.51%  '
$ws.Range("E51").Value = '  -32.31%  '
